# Fruta / hortaliza, semanal
# Insert a new weekly price record at row 348 for "Puerro" (Vega Modelo de
# Temuco), pushing the existing rows 348-362 down to 349-363.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 348; this shifts rows
# 348-362 down to 349-363 and keeps their data intact.
$ws.Rows.Item(348).Insert()

# Populate the new row 348 with the latest weekly observation.
$ws.Cells.Item(348, 1).Value = 10
$ws.Cells.Item(348, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(348, 3).Value = "La Araucanía"
$ws.Cells.Item(348, 4).Value = 45267
$ws.Cells.Item(348, 5).Value = 9
$ws.Cells.Item(348, 6).Value = 100112005
$ws.Cells.Item(348, 7).Value = "Puerro"
$ws.Cells.Item(348, 8).Value = "Azul de Maquehue"
$ws.Cells.Item(348, 9).Value = "Primera"
$ws.Cells.Item(348, 10).Value = 20
$ws.Cells.Item(348, 11).Value = 8000
$ws.Cells.Item(348, 12).Value = 8000
$ws.Cells.Item(348, 13).Value = 8000
$ws.Cells.Item(348, 14).Value = "`$/docena de paquetes"
$ws.Cells.Item(348, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(348, 16).Value = 667
$ws.Cells.Item(348, 17).Value = 12
$ws.Cells.Item(348, 18).Value = "Hortaliza"
